# Auto-generated edit script: adds observation rows 9-12 to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-BoolCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-TextCell($row, $col, $val) {
    # Leading apostrophe forces literal text (prevents Excel from
    # auto-converting date-like strings, e.g. "2023-09-14", into
    # date serial numbers). Resetting the style afterwards strips
    # the "quote prefix" formatting flag so the cell matches a
    # plain, unstyled text cell (same as the rest of the sheet).
    $ws.Cells.Item($row, $col).Value = "'" + $val
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# Row 9
Set-NumCell 9 1 112093593
Set-NumCell 9 2 103288
Set-TextCell 9 3 'Ovaliderad'
Set-TextCell 9 4 'LC'
Set-NumCell 9 5 221144
Set-TextCell 9 6 'Grönpyrola'
Set-TextCell 9 7 'Pyrola chlorantha'
Set-TextCell 9 8 'Sw.'
Set-TextCell 9 9 ''
Set-TextCell 9 11 'fullt utvecklade blad'
Set-TextCell 9 16 'Svarvartorp ca 400 m SO om, Upl'
Set-NumCell 9 17 653785.7777792643
Set-NumCell 9 18 6577035.071886262
Set-NumCell 9 19 10
Set-TextCell 9 20 'Stockholm'
Set-TextCell 9 21 'Ekerö'
Set-TextCell 9 22 'Uppland'
Set-TextCell 9 23 'Ekerö'
Set-TextCell 9 25 '2023-09-14'
Set-TextCell 9 26 '00:00'
Set-TextCell 9 27 '2023-09-14'
Set-TextCell 9 28 '00:00'
Set-BoolCell 9 30 $false
Set-BoolCell 9 31 $false
Set-BoolCell 9 33 $false
Set-TextCell 9 35 'Gles barrskog på sand (både tall och gran)'
Set-TextCell 9 46 ''
Set-TextCell 9 49 'Jan Yngve Andersson'
Set-TextCell 9 50 'Jan Yngve Andersson'
Set-TextCell 9 51 ''

# Row 10
Set-NumCell 10 1 112093595
Set-NumCell 10 2 90658
Set-TextCell 10 3 'Ovaliderad'
Set-TextCell 10 4 'NT'
Set-NumCell 10 5 4361
Set-TextCell 10 6 'Orange taggsvamp'
Set-TextCell 10 7 'Hydnellum aurantiacum'
Set-TextCell 10 8 '(Batsch:Fr.) P.Karst.'
Set-TextCell 10 9 ''
Set-TextCell 10 11 'teleomorf'
Set-TextCell 10 16 'Svarvartorp ca 400 m SO om, Upl'
Set-NumCell 10 17 653792.4227122802
Set-NumCell 10 18 6576997.511297328
Set-NumCell 10 19 10
Set-TextCell 10 20 'Stockholm'
Set-TextCell 10 21 'Ekerö'
Set-TextCell 10 22 'Uppland'
Set-TextCell 10 23 'Ekerö'
Set-TextCell 10 25 '2023-09-14'
Set-TextCell 10 26 '00:00'
Set-TextCell 10 27 '2023-09-14'
Set-TextCell 10 28 '00:00'
Set-BoolCell 10 30 $false
Set-BoolCell 10 31 $false
Set-BoolCell 10 33 $false
Set-TextCell 10 35 'Gles barrskog på sand (både tall och gran)'
Set-TextCell 10 46 ''
Set-TextCell 10 49 'Jan Yngve Andersson'
Set-TextCell 10 50 'Jan Yngve Andersson'
Set-TextCell 10 51 ''

# Row 11
Set-NumCell 11 1 112093592
Set-NumCell 11 2 90658
Set-TextCell 11 3 'Ovaliderad'
Set-TextCell 11 4 'NT'
Set-NumCell 11 5 4361
Set-TextCell 11 6 'Orange taggsvamp'
Set-TextCell 11 7 'Hydnellum aurantiacum'
Set-TextCell 11 8 '(Batsch:Fr.) P.Karst.'
Set-TextCell 11 9 ''
Set-TextCell 11 11 'teleomorf'
Set-TextCell 11 16 'Svarvartorp ca 400 m SO om, Upl'
Set-NumCell 11 17 653789.0938753984
Set-NumCell 11 18 6577029.07181866
Set-NumCell 11 19 10
Set-TextCell 11 20 'Stockholm'
Set-TextCell 11 21 'Ekerö'
Set-TextCell 11 22 'Uppland'
Set-TextCell 11 23 'Ekerö'
Set-TextCell 11 25 '2023-09-14'
Set-TextCell 11 26 '00:00'
Set-TextCell 11 27 '2023-09-14'
Set-TextCell 11 28 '00:00'
Set-BoolCell 11 30 $false
Set-BoolCell 11 31 $false
Set-BoolCell 11 33 $false
Set-TextCell 11 35 'Gles barrskog på sand (både tall och gran)'
Set-TextCell 11 46 ''
Set-TextCell 11 49 'Jan Yngve Andersson'
Set-TextCell 11 50 'Jan Yngve Andersson'
Set-TextCell 11 51 ''

# Row 12
Set-NumCell 12 1 112097135
Set-NumCell 12 2 90666
Set-TextCell 12 3 'Ovaliderad'
Set-TextCell 12 4 'LC'
Set-NumCell 12 5 4364
Set-TextCell 12 6 'Dropptaggsvamp'
Set-TextCell 12 7 'Hydnellum ferrugineum'
Set-TextCell 12 8 '(Fr.:Fr.) P. Karst.'
Set-TextCell 12 9 ''
Set-TextCell 12 10 ''
Set-TextCell 12 11 'teleomorf'
Set-TextCell 12 14 ''
Set-TextCell 12 16 'Svarvartorp ca 400 m SO om, Upl'
Set-NumCell 12 17 653888.520037169
Set-NumCell 12 18 6576888.42942148
Set-NumCell 12 19 10
Set-TextCell 12 20 'Stockholm'
Set-TextCell 12 21 'Ekerö'
Set-TextCell 12 22 'Uppland'
Set-TextCell 12 23 'Ekerö'
Set-TextCell 12 25 '2023-09-14'
Set-TextCell 12 26 '00:00'
Set-TextCell 12 27 '2023-09-14'
Set-TextCell 12 28 '00:00'
Set-TextCell 12 29 'Mörkröda droppar på hattöversidan. Smak besk efter ett långt tag, ej brännande. Köttet färgas mörkviolett med KOH.'
Set-BoolCell 12 30 $false
Set-BoolCell 12 31 $false
Set-TextCell 12 32 ''
Set-BoolCell 12 33 $false
Set-TextCell 12 35 'Gles barrskog på sand (både tall och gran)'
Set-TextCell 12 46 ''
Set-TextCell 12 49 'Jan Yngve Andersson'
Set-TextCell 12 50 'Jan Yngve Andersson'
Set-TextCell 12 51 ''

Write-Output "Added rows 9-12"
